# "Change planning excel file to reassign sah binning"
#
# Row 20 of the "workload" sheet is the task
# "SAH+binning as splitting criterion for BVH" (column B, shared string #10).
# Its time-allocation percentage is moved from member 3 (column F) to
# member 2 (column E): E20 goes 0 -> 100, F20 goes 100 -> 0.
# Downstream SUMPRODUCT/total formulas in rows 25-26 recompute automatically.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("workload")
$ws.Activate()

$ws.Range("E20").Value = 100
$ws.Range("F20").Value = 0

# Reflect the view state (scroll position / active cell) at save time.
$excel.ActiveWindow.ScrollRow = 2
$excel.ActiveWindow.ScrollColumn = 1
$ws.Range("C14").Select()
